# Regenerate merged AHB files
#
# The source data for this worksheet's "group header" rows (one per SG2
# segment group: rows 23, 27, 31, 36, 43, 47, 54) had not yet picked up the
# shared formatting used elsewhere in the table (gray fill + border, no
# "ÄNDERUNG" marker in the L/"Änderung" column). The remaining data rows in
# each group only needed their L column's stray "ÄNDERUNG" marker cleared.
#
# We reuse row 2 (already in the correct, final "group header" format) as
# the formatting template for the header rows, and its L cell (already
# blank, correctly styled) as the template for clearing the "Änderung"
# column on the rest of the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group header rows: copy the full A:V formatting from the reference
#     header row (row 2). (Note: ClearContents() drops CutCopyMode, so all
#     the paste operations for a given clipboard source are done first,
#     and the content clears happen afterwards in a separate pass.)
$headerRows = @(23, 27, 31, 36, 43, 47, 54)
$ws.Range("A2:V2").Copy() | Out-Null
foreach ($r in $headerRows) {
    $ws.Range("A$r`:V$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Remaining rows in each group: copy the (already blank) L2 cell's
#     format onto each row's L cell.
$lRows = @(24,25,26, 28,29,30, 32,33,34, 37,38,39,40,41,42, 44,45,46, 49,50,51,52,53, 55,56,57,58,59,60)
$ws.Range("L2").Copy() | Out-Null
foreach ($r in $lRows) {
    $ws.Range("L$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Finally, clear the stray "ÄNDERUNG" content out of every L cell
#     touched above (header rows included) now that formats are in place.
$allLRows = $headerRows + $lRows
foreach ($r in $allLRows) {
    $ws.Range("L$r").ClearContents() | Out-Null
}

Write-Host "Applied merged-AHB regeneration formatting fixups to rows 23-60."
